$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 6) appended to the inquiries table.
$ws.Range("A6").Value = "2025-11-14T04:32:20.098Z"
$ws.Range("B6").Value = "de"
$ws.Range("C6").Value = "eee"
$ws.Range("D6").Value = "dev@gmail.com"

# Phone number: keep as text (matches existing text-typed phone column)
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3142977875"

$ws.Range("F6").Value = "employer"
$ws.Range("G6").Value = "efhf"

# Student ID / Current Company are blank for this submission.
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""

$ws.Range("J6").Value = "Yes"
$ws.Range("K6").Value = 3

# Upcoming Event id: keep as text (matches existing text-typed column)
$ws.Range("L6").NumberFormat = "@"
$ws.Range("L6").Value = "1140"

$ws.Range("M6").Value = "def"
